# Scheduled runner update: refresh market-board derived figures
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) for a
# handful of Leve rows across the crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2971.3333
$ws.Range("I80").Value = 925.0625
$ws.Range("J80").Value = 7063.875
$ws.Range("K80").Value = 2775.1875
$ws.Range("L80").Value = 21191.625
$ws.Range("M80").Value = -1777.1875
$ws.Range("N80").Value = -23187.625

$ws.Range("H83").Value = 2971.3333
$ws.Range("I83").Value = 925.0625
$ws.Range("J83").Value = 7063.875
$ws.Range("K83").Value = 8325.5625
$ws.Range("L83").Value = 63574.875
$ws.Range("M83").Value = -3333.5625
$ws.Range("N83").Value = -73558.875

$ws.Range("H99").Value = 735.2308
$ws.Range("I99").Value = 590.8
$ws.Range("J99").Value = 1216.6666
$ws.Range("K99").Value = 1772.4
$ws.Range("L99").Value = 3649.9998
$ws.Range("M99").Value = -274.3999999999999
$ws.Range("N99").Value = -6645.9998

$ws.Range("H129").Value = 993.8393
$ws.Range("I129").Value = 251.57143
$ws.Range("J129").Value = 1099.8776
$ws.Range("K129").Value = 754.71429
$ws.Range("L129").Value = 3299.6328
$ws.Range("M129").Value = 4245.28571
$ws.Range("N129").Value = -13299.6328

$ws.Range("H134").Value = 43400
$ws.Range("J134").Value = 43400
$ws.Range("L134").Value = 43400
$ws.Range("N134").Value = -53540

$ws.Range("H137").Value = 1029.0256
$ws.Range("I137").Value = 758.96155
$ws.Range("J137").Value = 1569.1538
$ws.Range("K137").Value = 2276.88465
$ws.Range("L137").Value = 4707.4614
$ws.Range("M137").Value = 273.11535
$ws.Range("N137").Value = -9807.4614

$ws.Range("H139").Value = 70180
$ws.Range("J139").Value = 70180
$ws.Range("L139").Value = 70180
$ws.Range("N139").Value = -80460

$ws.Range("H140").Value = 70840.91
$ws.Range("J140").Value = 89906.25
$ws.Range("L140").Value = 89906.25
$ws.Range("N140").Value = -100266.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 111117390
$ws.Range("J37").Value = 8099
$ws.Range("L37").Value = 8099
$ws.Range("N37").Value = -8645

$ws.Range("H74").Value = 906.45
$ws.Range("I74").Value = 616.4286
$ws.Range("J74").Value = 1583.1666
$ws.Range("K74").Value = 616.4286
$ws.Range("L74").Value = 1583.1666
$ws.Range("M74").Value = 257.5714
$ws.Range("N74").Value = -3331.1666

$ws.Range("H77").Value = 906.45
$ws.Range("I77").Value = 616.4286
$ws.Range("J77").Value = 1583.1666
$ws.Range("K77").Value = 3082.143
$ws.Range("L77").Value = 7915.833000000001
$ws.Range("M77").Value = 1285.857
$ws.Range("N77").Value = -16651.833

$ws.Range("H122").Value = 5084.724
$ws.Range("I122").Value = 5303.4546
$ws.Range("J122").Value = 4397.2856
$ws.Range("K122").Value = 15910.3638
$ws.Range("L122").Value = 13191.8568
$ws.Range("M122").Value = -13460.3638
$ws.Range("N122").Value = -18091.8568

$ws.Range("H138").Value = 67712.5
$ws.Range("J138").Value = 67712.5
$ws.Range("L138").Value = 67712.5
$ws.Range("N138").Value = -77992.5

$ws.Range("H139").Value = 60833.332
$ws.Range("J139").Value = 60833.332
$ws.Range("L139").Value = 60833.332
$ws.Range("N139").Value = -71113.33199999999

$ws.Range("H140").Value = 90575
$ws.Range("J140").Value = 90575
$ws.Range("L140").Value = 90575
$ws.Range("N140").Value = -100935

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 65140
$ws.Range("J140").Value = 65140
$ws.Range("L140").Value = 65140
$ws.Range("N140").Value = -75500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4557.9067
$ws.Range("I31").Value = 3413.25
$ws.Range("J31").Value = 5236.222
$ws.Range("K31").Value = 3413.25
$ws.Range("L31").Value = 5236.222
$ws.Range("M31").Value = -3118.25
$ws.Range("N31").Value = -5826.222

$ws.Range("H34").Value = 4557.9067
$ws.Range("I34").Value = 3413.25
$ws.Range("J34").Value = 5236.222
$ws.Range("K34").Value = 3413.25
$ws.Range("L34").Value = 5236.222
$ws.Range("M34").Value = -3211.25
$ws.Range("N34").Value = -5640.222

$ws.Range("H50").Value = 7581.8335
$ws.Range("J50").Value = 7581.8335
$ws.Range("L50").Value = 7581.8335
$ws.Range("N50").Value = -8831.833500000001

$ws.Range("H51").Value = 8219.4
$ws.Range("J51").Value = 9274.25
$ws.Range("L51").Value = 9274.25
$ws.Range("N51").Value = -10746.25

$ws.Range("H60").Value = 7923.75
$ws.Range("J60").Value = 8200.666999999999
$ws.Range("L60").Value = 8200.666999999999
$ws.Range("N60").Value = -9222.666999999999

$ws.Range("H61").Value = 8219.4
$ws.Range("J61").Value = 9274.25
$ws.Range("L61").Value = 9274.25
$ws.Range("N61").Value = -9970.25

$ws.Range("H68").Value = 15824.25
$ws.Range("J68").Value = 15824.25
$ws.Range("L68").Value = 15824.25
$ws.Range("N68").Value = -17322.25

$ws.Range("H71").Value = 15824.25
$ws.Range("J71").Value = 15824.25
$ws.Range("L71").Value = 47472.75
$ws.Range("N71").Value = -54960.75

$ws.Range("H138").Value = 47300
$ws.Range("J138").Value = 47300
$ws.Range("L138").Value = 47300
$ws.Range("N138").Value = -57580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 72231.625
$ws.Range("I80").Value = 103536.91
$ws.Range("K80").Value = 103536.91
$ws.Range("M80").Value = -102538.91

$ws.Range("H83").Value = 72231.625
$ws.Range("I83").Value = 103536.91
$ws.Range("K83").Value = 517684.55
$ws.Range("M83").Value = -512692.55

$ws.Range("H140").Value = 89974.5
$ws.Range("J140").Value = 89974.5
$ws.Range("L140").Value = 89974.5
$ws.Range("N140").Value = -100334.5

$ws.Range("H141").Value = 67950
$ws.Range("J141").Value = 67950
$ws.Range("L141").Value = 67950
$ws.Range("N141").Value = -78310

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3622.3076
$ws.Range("I93").Value = 4066.5
$ws.Range("J93").Value = 3424.889
$ws.Range("K93").Value = 4066.5
$ws.Range("L93").Value = 3424.889
$ws.Range("M93").Value = -2818.5
$ws.Range("N93").Value = -5920.889

$ws.Range("H136").Value = 3558.8538
$ws.Range("I136").Value = 3107
$ws.Range("J136").Value = 3793.1482
$ws.Range("K136").Value = 9321
$ws.Range("L136").Value = 11379.4446
$ws.Range("M136").Value = -6771
$ws.Range("N136").Value = -16479.4446

$ws.Range("H140").Value = 58865
$ws.Range("J140").Value = 59166.668
$ws.Range("L140").Value = 59166.668
$ws.Range("N140").Value = -69526.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 34888000
$ws.Range("I132").Value = 46876624
$ws.Range("J132").Value = 12000.728
$ws.Range("K132").Value = 140629872
$ws.Range("L132").Value = 36002.18399999999
$ws.Range("M132").Value = -140627342
$ws.Range("N132").Value = -41062.18399999999

$ws.Range("H136").Value = 1020.5417
$ws.Range("I136").Value = 885.5909
$ws.Range("J136").Value = 2505
$ws.Range("K136").Value = 2656.7727
$ws.Range("L136").Value = 7515
$ws.Range("M136").Value = -106.7727
$ws.Range("N136").Value = -12615

$ws.Range("H138").Value = 57418.43
$ws.Range("J138").Value = 57418.43
$ws.Range("L138").Value = 57418.43
$ws.Range("N138").Value = -67698.42999999999

$ws.Range("H139").Value = 61840
$ws.Range("J139").Value = 61840
$ws.Range("L139").Value = 61840
$ws.Range("N139").Value = -72120
